# Update cryptos list with latest prices / volume changes, and re-sort a
# few coins whose relative ranking changed (Chainlink/Uniswap swap,
# BabyDogeCoin inserted ahead of Cronos/Algorand, USDD replaced by Algorand).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values that look numeric (e.g. "19.20",
# "0.0516") but must stay text so formatting (trailing zeros, thousands
# dots, subscript notation) is preserved exactly as scraped. Force the
# whole column to text format before writing any values.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.793.86"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.640.80"
$ws.Range("E3").Value = "  -0.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.20%  "

# Row 5 - BNB
$ws.Range("D5").Value = "218.56"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.07%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.29%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.10%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.79%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.20"
$ws.Range("E10").Value = "  +0.20%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.53%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.869.91"
$ws.Range("E12").Value = "  -0.18%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.642.29"
$ws.Range("E13").Value = "  -0.04%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -0.67%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  -0.39%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "65.01"
$ws.Range("E16").Value = "  +0.80%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.810.01"
$ws.Range("E17").Value = "  +0.00%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -0.86%  "

# Row 19 - BitcoinCash
$ws.Range("E19").Value = "  +1.07%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.32%  "

# Row 21 - was Chainlink, now Uniswap
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "4.36"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22 - was Uniswap, now Chainlink
$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  +4.60%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "2.36"
$ws.Range("E23").Value = "  -2.01%  "

# Row 24 - Avalanche
$ws.Range("D24").Value = "9.15"

# Row 25 - Monero
$ws.Range("D25").Value = "147.60"
$ws.Range("E25").Value = "  +1.67%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.29%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.25%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "7.09"
$ws.Range("E28").Value = "  +0.28%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "15.70"
$ws.Range("E29").Value = "  +0.06%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -1.26%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +1.01%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +1.89%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -0.85%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +0.83%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.264.31"
$ws.Range("E35").Value = "  -2.09%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.27%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.0174"
$ws.Range("E37").Value = "  -0.51%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -1.50%  "

# Row 39 - ARBITRUM
$ws.Range("D39").Value = "0.816"
$ws.Range("E39").Value = "  -1.19%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  -0.21%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  -0.59%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "5.34"
$ws.Range("E42").Value = "  -0.43%  "

# Row 43 - RocketPoolETH
$ws.Range("D43").Value = "1.780.25"
$ws.Range("E43").Value = "  -0.84%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  -4.19%  "

# Row 45 - Quant
$ws.Range("D45").Value = "92.39"
$ws.Range("E45").Value = "  +1.13%  "

# Row 46 - Aave
$ws.Range("D46").Value = "60.69"
$ws.Range("E46").Value = "  +0.28%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "1.58"
$ws.Range("E47").Value = "  -0.98%  "

# Row 48 - was Cronos, now BabyDogeCoin
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -0.38%  "

# Row 49 - was Algorand, now Cronos
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0516"
$ws.Range("E49").Value = "  -0.59%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "7.56"
$ws.Range("E50").Value = "  -1.97%  "

# Row 51 - was USDD, now Algorand
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0962"
$ws.Range("E51").Value = "  -1.65%  "
